$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shorten two of the "note" explanations (column H) ---
# Edited in this order (Wars of the Roses first, then Mona Lisa) to match
# the shared-string append order of the source edit.

# Row 65: "ばら戦争が起きた場所は？" -> note column H
$ws.Cells.Item(65, 8).Value = "薔薇戦争（ばらせんそう、英: Wars of the Roses）は、百年戦争終戦後に発生したイングランド中世封建諸侯による内乱"

# Row 63: "絵画、モナ・リザはどちらの手が上にあるか？" -> note column H
$ws.Cells.Item(63, 8).Value = "1503年 - 1519年頃、レオナルド・ダヴィンチによって描かれた絵画。もともとはフランス王フランソワ1世が購入した作品だが、現在はフランスの国有財産"

# --- Extend the row-numbering formula in column A down through row 66 ---
# Originally the shared formula only covered A4:A55; rows 56-66 had no
# index number yet. Fill them in with the same "=A(row-1)+1" pattern.
for ($r = 56; $r -le 66; $r++) {
    $ws.Cells.Item($r, 1).Formula = "=A" + ($r - 1) + "+1"
}

# --- Update the view: scroll/select so column A's new fill is visible ---
$ws.Activate() | Out-Null
$ws.Range("A55:A66").Select() | Out-Null
